$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Forward (Cold) Post ydist -> "Forward (Cold) Post ydist far"
$ws.Range("A1").Value = "Forward (Cold) Post ydist far"

# Row 2: Aft (Hot )Post ydist -> "Aft (Hot )Post ydist far"
$ws.Range("A2").Value = "Aft (Hot )Post ydist far"

# Rows 3 & 4 unchanged (Forward/Aft xdist)

# Row 5 / Row 6 labels unchanged in content (Engine Centerpoint ydist/xdist), just reordered string ids - no value change needed
# Row 7: MaxThrust value changes from 71220 to 45500 (style/border index for A7/B7 stays as-is)
$ws.Range("B7").Value = 45500

# New Row 8: Forward (Cold) Post ydist close
$ws.Range("A8").Value = "Forward (Cold) Post ydist close"
$ws.Range("B8").Value = 2.68

# New Row 9: Aft (Hot )Post ydist close
$ws.Range("A9").Value = "Aft (Hot )Post ydist close"
$ws.Range("B9").Value = 2.68

# Copy formatting (style) from row 1 (A1:B1) to new rows 8 & 9 (A/B columns) since they use style s="1"/s="2" like row1/2
$ws.Range("A1:B1").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)  # xlPasteFormats = -4122
$ws.Range("A1:B1").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)

# Restore the values since PasteSpecial formats only shouldn't touch values, but ensure correctness anyway
$ws.Range("A8").Value = "Forward (Cold) Post ydist close"
$ws.Range("B8").Value = 2.68
$ws.Range("A9").Value = "Aft (Hot )Post ydist close"
$ws.Range("B9").Value = 2.68
